# "edições para trabalhar no caos"
# Anonymize the survey headers into short codes (ID, P1..P9) and fill in
# the previously-blank "external student" answer with "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Row 1: replace the long Portuguese questions with short column codes ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "P1"
$ws.Range("C1").Value = "P2"
$ws.Range("D1").Value = "P3"
$ws.Range("E1").Value = "P4"
$ws.Range("F1").Value = "P5"
$ws.Range("G1").Value = "P6"
$ws.Range("H1").Value = "P7"
$ws.Range("I1").Value = "P8"
$ws.Range("J1").Value = "P9"

# --- Row 2: the external-program answer was blank; mark it as not applicable ---
$ws.Range("E2").Value = "NA"

# --- Give every header cell the same look (centered, no border) ---
$hdr = $ws.Range("A1:J1")
$hdr.Font.Name = "Calibri"
$hdr.Font.Size = 11
$hdr.Font.ThemeColor = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.Borders.LineStyle = -4142

# --- Shrink the columns now that the headers are short codes ---
$ws.Range("A:A").ColumnWidth = 1.9986979166666665
$ws.Range("B:B").ColumnWidth = 7.666666666666667
$ws.Range("C:C").ColumnWidth = 24.166666666666668
$ws.Range("D:D").ColumnWidth = 3.1666666666666665
$ws.Range("E:E").ColumnWidth = 2.3307291666666665
$ws.Range("F:F").ColumnWidth = 32.998697916666664
$ws.Range("G:G").ColumnWidth = 11.166666666666666
$ws.Range("H:H").ColumnWidth = 33.166666666666664
$ws.Range("I:I").ColumnWidth = 3.4986979166666665
$ws.Range("J:J").ColumnWidth = 29.166666666666668

# --- Match the new selection left by the edit ---
$ws.Range("A1:J1").Select()
